$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 83.59375
$ws.Range("C2").Value = 83.59375
$ws.Range("D2").Value = 83.59375
$ws.Range("E2").Value = 66.66666666666667
$ws.Range("F2").Value = 57.55208333333333
$ws.Range("G2").Value = 54.42708333333333
$ws.Range("H2").Value = 53.90625
$ws.Range("I2").Value = 53.90625
